$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new column before column N (shifts N->O, O->P, P->Q)
$ws.Columns("N").Insert()

# Select O8 as the final selection (matches target selection state)
$ws.Range("O8").Select()
